# Assignment 4 working commit
# - Insert a "Remain" (Amt - Paid) column between Paid (E) and the old Int
#   column, pushing Int/Recurring from F/G to G/H.
# - Fill the new Remain column with D-E formulas (row 3 standalone, rows
#   4-11 as one fill/shared formula).
# - Backfill a few missing "Paid" amounts.
# - Turn the Us Bank "Paid" cell into a formula (1950 + interest).
# - Add a totals row (12) summing the Amt and Remain columns.
# - Leave the selection on E11, matching the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F. Existing F (Int) -> G, G (Recurring) -> H.
$ws.Columns("F:F").Insert()

# F3 = Amt - Paid (kept separate from the fill below, like the source file)
$ws.Range("F3").Formula = "=D3-E3"

# Fill in previously-blank "Paid" values
$ws.Range("E4").Value = 180.98
$ws.Range("E7").Value = 208.17
$ws.Range("E11").Value = 1068.9000000000001

# F4:F11 = Amt - Paid, filled as one formula across the range
$ws.Range("F4:F11").Formula = "=D4-E4"

# Us Bank's "Paid" amount becomes a formula: 1950 plus the interest charge
$ws.Range("E10").Formula = "=1950+C10"

# New totals row
$ws.Range("D12").Formula = "=SUM(D3:D11)"
$ws.Range("F12").Formula = "=SUM(F3:F11)"

# Match the saved selection
$ws.Range("E11").Select() | Out-Null
